# "updates for the new week"
# Insert a new week's row above the current last-entry (highlighted) row.
# Because Excel inherits row formatting from the row directly above an
# insertion point, this naturally "promotes" the previous week (row 26,
# unhighlighted) onto the old highlighted row's position (now row 27) and
# leaves the original highlighted row (with its yellow fill) pushed down
# to become the new current-week row (row 28), which is exactly the
# visual effect the author wants: last week's entry loses the highlight
# and this week's entry becomes the new highlighted "latest" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the old row 27 (1/20/2020, highlighted, 15 hrs) down to row 28,
# inserting a fresh (unhighlighted, date-formatted) row 27 above it.
$ws.Rows.Item(27).Insert()

# Row 27 now holds last week's corrected entry (no longer the latest).
$weekOf0120 = Get-Date -Year 2020 -Month 1 -Day 20 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("A27").Value = $weekOf0120
$ws.Range("B27").Value = 3

# Row 28 becomes this week's new entry (keeps the inherited highlighted style).
$weekOf0127 = Get-Date -Year 2020 -Month 1 -Day 27 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("A28").Value = $weekOf0127
$ws.Range("B28").Value = 15

# Leave the active selection on the newly entered hours cell.
$ws.Range("B28").Select()
